# edit.ps1 -- applies the "Añado titulo al archivo don-quijote" change.
#
# Summary of the edit (see xml_diff):
#   1. Insert a brand-new first paragraph containing the bold, bigger
#      "TITULO"/"TÍTULO" heading (Roboto Condensed, 20pt, bold, purple),
#      carrying the document's "_GoBack" bookmark.
#   2. The poem paragraphs below are otherwise unchanged (they simply
#      shift down by one position).
#   3. The closing "`git commit -m ...`" paragraph loses the "_GoBack"
#      bookmark (it moved up to the new heading paragraph) and gets its
#      "...la M" / "ancha"`" split redrawn as "...la " / "Mancha"`",
#      now wrapped in a grammar-check (gramStart/gramEnd) proofErr pair
#      instead of the bookmark that used to sit there.
#
# Because the COM Font/Bold/BoldBi properties don't let us reproduce the
# exact <w:b/><w:bCs/> / <w:sz/><w:szCs/> pairing on both the run AND the
# paragraph-mark (pPr/rPr) in one shot, both paragraphs are rewritten with
# Range.InsertXML using literal WordprocessingML -- this is just "paste
# this formatted run/paragraph" from the OM's point of view and keeps
# every other part of the package (styles, fonts, settings, sectPr, the
# untouched middle paragraphs, ...) exactly as Word would leave them.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert the new "TÍTULO" heading paragraph before everything else.
# ---------------------------------------------------------------------
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Roboto Condensed" w:hAnsi="Roboto Condensed"/><w:b/><w:bCs/><w:color w:val="7030A0"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed" w:hAnsi="Roboto Condensed"/><w:b/><w:bCs/><w:color w:val="7030A0"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>TÍTULO</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2. Rewrite the closing backtick "`git commit -m ...`" paragraph so the
#    "_GoBack" bookmark is gone and "Mancha" is a single gram-checked run.
# ---------------------------------------------------------------------
$closingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t>`</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve"> -m “Don Quijote de la </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="7030A0"/></w:rPr><w:t>Mancha”`</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$closingParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingParagraph.Range.InsertXML($closingXml)
